# Updating Human Sensing parts - add IR sensor connectors/housings
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 28: Thermal Sensor Connectors ---
$ws.Range("B28").Value = "Thermal Sensor Connectors"
$ws.Range("C28").Value = "SSHL-002T-P0.2"
$ws.Range("D28").Value = 10
$ws.Range("E28").Value = 0.105
$ws.Range("G28").Value = "http://www.digikey.com/product-detail/en/SSHL-002T-P0.2/455-1606-1-ND/1642795"
$ws.Range("H28").Value = "Connector for IR"
$ws.Range("I28").Value = "10 includes replacements"

# --- Row 29: Thermal Sensor Housing ---
$ws.Range("B29").Value = "Thermal Sensor Housing"
$ws.Range("C29").Value = "GHR-04V-S"
$ws.Range("D29").Value = 8
$ws.Range("E29").Value = 0.13
$ws.Range("G29").Value = "http://www.digikey.com/product-detail/en/GHR-04V-S/455-1594-ND/807816"
$ws.Range("H29").Value = "Housing for IR"
$ws.Range("I29").Value = "8 Includes replacements"

# Hyperlinks for the new parts (digikey product pages)
$ws.Hyperlinks.Add($ws.Range("G28"), "http://www.digikey.com/product-detail/en/SSHL-002T-P0.2/455-1606-1-ND/1642795")
$ws.Hyperlinks.Add($ws.Range("G29"), "http://www.digikey.com/product-detail/en/GHR-04V-S/455-1594-ND/807816")

# Restore the standard "hyperlink" cell style (matches other link cells like G5/G27)
# without leaving behind Hyperlinks.Add's auto-generated style on the cell itself.
$ws.Range("G5").Copy()
$ws.Range("G28").PasteSpecial(-4122)
$ws.Range("G5").Copy()
$ws.Range("G29").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- View state: zoom + selection as left by the editor ---
$ws.Activate()
$excel.ActiveWindow.Zoom = 80
$ws.Range("E30").Select()
